$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new Argent (Solar) price row for 2025-05-17, mirroring the
# existing data rows. All the source values are stored as plain text
# (not numbers/dates), so force the target range to Text format before
# assigning, then restore the "Normal" style so no extraneous number
# formatting lingers on the new cells.
$newRow = 77
$rngAddr = "A" + $newRow + ":J" + $newRow
$rng = $ws.Range($rngAddr)
$rng.NumberFormat = "@"

$ws.Range("A$newRow").Value = "2025-05-17"
$ws.Range("B$newRow").Value = "37.5"
$ws.Range("C$newRow").Value = "37"
$ws.Range("D$newRow").Value = "0.94"
$ws.Range("E$newRow").Value = "0.258"
$ws.Range("F$newRow").Value = "0.09"
$ws.Range("G$newRow").Value = "5,298"
$ws.Range("H$newRow").Value = "7,931"
$ws.Range("I$newRow").Value = "7,981"
$ws.Range("J$newRow").Value = "7.2226"

$rng.Style = "Normal"
